$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29 (WIN row): the "WIN" popup text changes from the old
# congratulations/win copy to a generic "END GAME!" / "GAME KẾT THÚC!" pair
# (the previous strings are no longer referenced anywhere, so the shared
# string table entries are dropped automatically).
$ws.Range("B29").Value = "END GAME!"
$ws.Range("C29").Value = "GAME KẾT THÚC!"

# New row 36: first "players set" / character-set entry, added at the end
# of the table (extends the used range from A1:C35 to A1:C36).
$ws.Range("A36").Value = "CHARSET"
$ws.Range("B36").Value = "CHOOSE PLAYERS SET: "
$ws.Range("C36").Value = "CHỌN SET NHÂN VẬT: "

# Update the view selection to match the saved workbook state.
$ws.Range("C15").Select()
